$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.095.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.05%  '

$ws.Range("D3").Value = "'1.836.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.15%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.41%  '

$ws.Range("D5").Value = "'242.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.97%  '

$ws.Range("D6").Value = "'0.6161"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.08%  '

$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("D8").Value = "'0.07472"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.14%  '

$ws.Range("D9").Value = "'0.2929"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.72%  '

$ws.Range("D10").Value = "'23.09"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.97%  '

$ws.Range("D11").Value = "'0.07695"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.42%  '

$ws.Range("D12").Value = "'1.827.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.38%  '

$ws.Range("D13").Value = "'4.995"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.33%  '

$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").Value = "'82.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.73%  '

$ws.Range("D16").Value = "'0.000009133"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.11%  '

$ws.Range("D17").Value = "'5.915"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.02%  '

$ws.Range("D18").Value = "'29.066.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.15%  '

$ws.Range("D19").Value = "'2.082.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.25%  '

$ws.Range("D20").Value = "'232.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.11%  '

$ws.Range("D21").Value = "'12.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.16%  '

$ws.Range("D22").Value = "'1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.37%  '

$ws.Range("D23").Value = "'7.189"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '

$ws.Range("D24").Value = "'1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.42%  '

$ws.Range("D25").Value = "'159.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.90%  '

$ws.Range("D26").Value = "'0.1401"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.80%  '

$ws.Range("D27").Value = "'8.491"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.80%  '

$ws.Range("D28").Value = "'17.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.80%  '

$ws.Range("D29").Value = "'1.500"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").Value = "'4.155"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("D31").Value = "'4.107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.66%  '

$ws.Range("D32").Value = "'0.05503"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.53%  '

$ws.Range("D33").Value = "'1.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("D34").Value = "'1.834"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.62%  '

$ws.Range("D35").Value = "'0.7372"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.47%  '

$ws.Range("D36").Value = "'1.140"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '

$ws.Range("D37").Value = "'2.663"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.03%  '

$ws.Range("D38").Value = "'2.775"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.60%  '

$ws.Range("D39").Value = "'0.01780"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.58%  '

$ws.Range("D40").Value = "'1.213.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.96%  '

$ws.Range("D41").Value = "'6.443"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.46%  '

$ws.Range("D42").Value = "'0.8935"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.17%  '

$ws.Range("D43").Value = "'1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").Value = "'101.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.18%  '

$ws.Range("D45").Value = "'1.980.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("D47").Value = "'0.5095"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("D48").Value = "'0.00000000118"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.88%  '

$ws.Range("D49").Value = "'0.4074"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.24%  '

$ws.Range("D50").Value = "'9.109"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.10%  '

$ws.Range("D51").Value = "'0.05824"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.53%  '
